$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.2143979625585434
$ws.Range("C2").Value = 1.424134207528472
$ws.Range("D2").Value = 7.525940274036349
$ws.Range("E2").Value = 2.743344723879292
$ws.Range("F2").Value = 2.796421391031461
$ws.Range("G2").Value = 23

# Row 3 (Q1)
$ws.Range("B3").Value = 0.4670168130507297
$ws.Range("C3").Value = 1.853629581901463
$ws.Range("D3").Value = 13.55051072646238
$ws.Range("E3").Value = 3.681101835926627
$ws.Range("F3").Value = 3.737282769615654
$ws.Range("G3").Value = 22

# Row 4 (Q2)
$ws.Range("B4").Value = -0.8354076150457171
$ws.Range("C4").Value = 1.386534768542286
$ws.Range("D4").Value = 7.40374990259892
$ws.Range("E4").Value = 2.720983260257019
$ws.Range("F4").Value = 2.653513938212625
$ws.Range("G4").Value = 21

# Row 5 (Q3)
$ws.Range("B5").Value = -0.04486636216424342
$ws.Range("C5").Value = 0.5141898899062147
$ws.Range("D5").Value = 0.5563440360261616
$ws.Range("E5").Value = 0.7458847337398464
$ws.Range("F5").Value = 0.7638758824313834
$ws.Range("G5").Value = 20

# Row 6 (Q4)
$ws.Range("B6").Value = 0.08415703371898316
$ws.Range("C6").Value = 0.7459901986017395
$ws.Range("D6").Value = 1.174992571369782
$ws.Range("E6").Value = 1.08397074285692
$ws.Range("F6").Value = 1.110312597020984
$ws.Range("G6").Value = 19

# Row 7 (Q5)
$ws.Range("B7").Value = -0.2038590712083605
$ws.Range("C7").Value = 0.6468456134312068
$ws.Range("D7").Value = 0.8625708541187559
$ws.Range("E7").Value = 0.9287469268421596
$ws.Range("F7").Value = 0.9323664388718182
$ws.Range("G7").Value = 18

# Row 8 (Q6)
$ws.Range("B8").Value = -0.1285815229469241
$ws.Range("C8").Value = 0.5264431630486278
$ws.Range("D8").Value = 0.5203200017941761
$ws.Range("E8").Value = 0.7213321022900452
$ws.Range("F8").Value = 0.7316238571562934
$ws.Range("G8").Value = 17

# Row 9 (Q7)
$ws.Range("B9").Value = 0.08851026789393576
$ws.Range("C9").Value = 0.4260801985527543
$ws.Range("D9").Value = 0.3054846069715224
$ws.Range("E9").Value = 0.5527066192579227
$ws.Range("F9").Value = 0.5634659783389385
$ws.Range("G9").Value = 16

# Row 10 (Q8)
$ws.Range("B10").Value = 0.007460085389865952
$ws.Range("C10").Value = 0.3714499647231649
$ws.Range("D10").Value = 0.2993920067773447
$ws.Range("E10").Value = 0.5471672566750908
$ws.Range("F10").Value = 0.5663192757087402
$ws.Range("G10").Value = 15

# Row 11 (Q9)
$ws.Range("B11").Value = 0.005873690370570317
$ws.Range("C11").Value = 0.3581201956056038
$ws.Range("D11").Value = 0.2817709375354674
$ws.Range("E11").Value = 0.5308210032915686
$ws.Range("F11").Value = 0.5524621921347255
